# تعديل تلقائي في شيت Card15 by admin at 2025-11-02 06:04:41
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card15")

# O1 header: "Serviced by " -> "Serviced by" (drop the trailing space)
$ws.Range("O1").Value = "Serviced by"

# M2: "nan" -> Arabic text "تم التشغيل"
$ws.Range("M2").Value = "تم التشغيل"

# O2:O12 were blank cells; give them the text "nan" like the rest of the
# "Correction"/other data columns on this row range.
for ($r = 2; $r -le 12; $r++) {
    $ws.Cells.Item($r, 15).Value = "nan"
}
